{"js": "// Replace each three-digit-by-one-digit multiplication prompt in the\n// worksheet table with its new equation. Every old value below is unique\n// in the document, so an exact, case-sensitive search-and-replace on the\n// whole body is safe and keeps each run's existing formatting (font,\n// size, etc.) intact.\nconst replacements = [\n  [\"865\u00d72=\", \"913\u00d72=\"],\n  [\"454\u00d75=\", \"661\u00d73=\"],\n  [\"570\u00d79=\", \"703\u00d73=\"],\n  [\"227\u00d79=\", \"358\u00d75=\"],\n  [\"698\u00d73=\", \"726\u00d79=\"],\n  [\"656\u00d73=\", \"124\u00d76=\"],\n  [\"301\u00d74=\", \"554\u00d78=\"],\n  [\"443\u00d79=\", \"683\u00d79=\"],\n  [\"955\u00d79=\", \"506\u00d79=\"],\n  [\"450\u00d76=\", \"785\u00d79=\"],\n  [\"459\u00d77=\", \"803\u00d78=\"],\n  [\"166\u00d77=\", \"180\u00d77=\"],\n  [\"139\u00d76=\", \"231\u00d73=\"],\n  [\"706\u00d73=\", \"152\u00d72=\"],\n  [\"974\u00d78=\", \"946\u00d76=\"],\n  [\"559\u00d72=\", \"392\u00d76=\"],\n  [\"545\u00d77=\", \"636\u00d79=\"],\n  [\"251\u00d75=\", \"611\u00d77=\"],\n  [\"251\u00d73=\", \"659\u00d74=\"],\n  [\"239\u00d78=\", \"667\u00d74=\"],\n  [\"391\u00d76=\", \"181\u00d77=\"],\n  [\"491\u00d76=\", \"961\u00d74=\"],\n  [\"814\u00d76=\", \"650\u00d74=\"],\n  [\"778\u00d79=\", \"328\u00d79=\"],\n  [\"910\u00d77=\", \"592\u00d72=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit-by-one-digit multiplication prompt in the\n# worksheet table with its new equation. Every old value is unique in the\n# document, so Find/Replace (one match each) is safe and keeps each run's\n# existing formatting (font, size, etc.) intact.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @{old = \"865\u00d72=\"; new = \"913\u00d72=\"},\n  @{old = \"454\u00d75=\"; new = \"661\u00d73=\"},\n  @{old = \"570\u00d79=\"; new = \"703\u00d73=\"},\n  @{old = \"227\u00d79=\"; new = \"358\u00d75=\"},\n  @{old = \"698\u00d73=\"; new = \"726\u00d79=\"},\n  @{old = \"656\u00d73=\"; new = \"124\u00d76=\"},\n  @{old = \"301\u00d74=\"; new = \"554\u00d78=\"},\n  @{old = \"443\u00d79=\"; new = \"683\u00d79=\"},\n  @{old = \"955\u00d79=\"; new = \"506\u00d79=\"},\n  @{old = \"450\u00d76=\"; new = \"785\u00d79=\"},\n  @{old = \"459\u00d77=\"; new = \"803\u00d78=\"},\n  @{old = \"166\u00d77=\"; new = \"180\u00d77=\"},\n  @{old = \"139\u00d76=\"; new = \"231\u00d73=\"},\n  @{old = \"706\u00d73=\"; new = \"152\u00d72=\"},\n  @{old = \"974\u00d78=\"; new = \"946\u00d76=\"},\n  @{old = \"559\u00d72=\"; new = \"392\u00d76=\"},\n  @{old = \"545\u00d77=\"; new = \"636\u00d79=\"},\n  @{old = \"251\u00d75=\"; new = \"611\u00d77=\"},\n  @{old = \"251\u00d73=\"; new = \"659\u00d74=\"},\n  @{old = \"239\u00d78=\"; new = \"667\u00d74=\"},\n  @{old = \"391\u00d76=\"; new = \"181\u00d77=\"},\n  @{old = \"491\u00d76=\"; new = \"961\u00d74=\"},\n  @{old = \"814\u00d76=\"; new = \"650\u00d74=\"},\n  @{old = \"778\u00d79=\"; new = \"328\u00d79=\"},\n  @{old = \"910\u00d77=\"; new = \"592\u00d72=\"}\n)\n\nforeach ($pair in $replacements) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  # wdFindContinue = 1, wdReplaceAll = 2\n  $find.Execute($pair.old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.new, 2) | Out-Null\n}\n"}
